$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of the "Fruta, Vega Modelo de Temuco - Coco" sheet:
# re-date each existing observation (column D) and update its
# Volumen / Precio minimo / Precio maximo / Precio promedio ponderado /
# Precio por Kg (columns M, N, O, P, S) to this weeks figures, then
# append one new observation as row 38.

# Row 2
$ws.Range("D2").Value = 44424
$ws.Range("M2").Value = 25
$ws.Range("N2").Value = 24000
$ws.Range("O2").Value = 24000
$ws.Range("P2").Value = 24000
$ws.Range("S2").Value = 1200

# Row 3
$ws.Range("D3").Value = 44214
$ws.Range("N3").Value = 25000
$ws.Range("O3").Value = 25000
$ws.Range("P3").Value = 25000
$ws.Range("S3").Value = 1250

# Row 4
$ws.Range("D4").Value = 44356
$ws.Range("M4").Value = 15

# Row 5
$ws.Range("D5").Value = 44349
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 24000
$ws.Range("O5").Value = 24000
$ws.Range("P5").Value = 24000
$ws.Range("S5").Value = 1200

# Row 6
$ws.Range("D6").Value = 44389
$ws.Range("M6").Value = 20

# Row 7
$ws.Range("D7").Value = 44390
$ws.Range("M7").Value = 10
$ws.Range("N7").Value = 24000
$ws.Range("O7").Value = 24000
$ws.Range("P7").Value = 24000
$ws.Range("S7").Value = 1200

# Row 8
$ws.Range("D8").Value = 44425

# Row 9
$ws.Range("D9").Value = 44232
$ws.Range("M9").Value = 15
$ws.Range("N9").Value = 25000
$ws.Range("O9").Value = 25000
$ws.Range("P9").Value = 25000
$ws.Range("S9").Value = 1250

# Row 10
$ws.Range("D10").Value = 44363
$ws.Range("M10").Value = 30

# Row 11
$ws.Range("D11").Value = 44431
$ws.Range("M11").Value = 40
$ws.Range("N11").Value = 24000
$ws.Range("O11").Value = 24000
$ws.Range("P11").Value = 24000
$ws.Range("S11").Value = 1200

# Row 12
$ws.Range("D12").Value = 44231
$ws.Range("N12").Value = 25000
$ws.Range("O12").Value = 25000
$ws.Range("P12").Value = 25000
$ws.Range("S12").Value = 1250

# Row 13
$ws.Range("D13").Value = 44418
$ws.Range("M13").Value = 20
$ws.Range("N13").Value = 24000
$ws.Range("O13").Value = 24000
$ws.Range("P13").Value = 24000
$ws.Range("S13").Value = 1200

# Row 14
$ws.Range("D14").Value = 44396
$ws.Range("M14").Value = 12

# Row 15
$ws.Range("D15").Value = 44221
$ws.Range("N15").Value = 25000
$ws.Range("O15").Value = 25000
$ws.Range("P15").Value = 25000
$ws.Range("S15").Value = 1250

# Row 16
$ws.Range("D16").Value = 44334
$ws.Range("N16").Value = 25000
$ws.Range("O16").Value = 25000
$ws.Range("P16").Value = 25000
$ws.Range("S16").Value = 1250

# Row 17
$ws.Range("D17").Value = 44419
$ws.Range("M17").Value = 40

# Row 18
$ws.Range("D18").Value = 44392
$ws.Range("M18").Value = 10

# Row 19
$ws.Range("D19").Value = 44249
$ws.Range("N19").Value = 25000
$ws.Range("O19").Value = 25000
$ws.Range("P19").Value = 25000
$ws.Range("S19").Value = 1250

# Row 20
$ws.Range("D20").Value = 44235
$ws.Range("M20").Value = 15
$ws.Range("N20").Value = 25000
$ws.Range("O20").Value = 25000
$ws.Range("P20").Value = 25000
$ws.Range("S20").Value = 1250

# Row 21
$ws.Range("D21").Value = 44433
$ws.Range("M21").Value = 10
$ws.Range("N21").Value = 24000
$ws.Range("O21").Value = 24000
$ws.Range("P21").Value = 24000
$ws.Range("S21").Value = 1200

# Row 22
$ws.Range("D22").Value = 44426
$ws.Range("N22").Value = 24000
$ws.Range("O22").Value = 24000
$ws.Range("P22").Value = 24000
$ws.Range("S22").Value = 1200

# Row 23
$ws.Range("D23").Value = 44421
$ws.Range("M23").Value = 20
$ws.Range("N23").Value = 24000
$ws.Range("O23").Value = 24000
$ws.Range("P23").Value = 24000
$ws.Range("S23").Value = 1200

# Row 24
$ws.Range("D24").Value = 44434
$ws.Range("M24").Value = 20

# Row 25
$ws.Range("D25").Value = 44222
$ws.Range("M25").Value = 15
$ws.Range("N25").Value = 25000
$ws.Range("O25").Value = 25000
$ws.Range("P25").Value = 25000
$ws.Range("S25").Value = 1250

# Row 26
$ws.Range("D26").Value = 44398

# Row 27
$ws.Range("D27").Value = 44420
$ws.Range("M27").Value = 35
$ws.Range("N27").Value = 25000
$ws.Range("O27").Value = 25000
$ws.Range("P27").Value = 25000
$ws.Range("S27").Value = 1250

# Row 28
$ws.Range("D28").Value = 44414
$ws.Range("M28").Value = 15
$ws.Range("N28").Value = 25000
$ws.Range("O28").Value = 25000
$ws.Range("P28").Value = 25000
$ws.Range("S28").Value = 1250

# Row 29
$ws.Range("D29").Value = 44432
$ws.Range("M29").Value = 30
$ws.Range("N29").Value = 24000
$ws.Range("O29").Value = 24000
$ws.Range("P29").Value = 24000
$ws.Range("S29").Value = 1200

# Row 30
$ws.Range("D30").Value = 44391
$ws.Range("M30").Value = 10
$ws.Range("N30").Value = 24000
$ws.Range("O30").Value = 24000
$ws.Range("P30").Value = 24000
$ws.Range("S30").Value = 1200

# Row 31
$ws.Range("D31").Value = 44251
$ws.Range("M31").Value = 15

# Row 32
$ws.Range("D32").Value = 44428
$ws.Range("M32").Value = 15

# Row 33
$ws.Range("D33").Value = 44442
$ws.Range("M33").Value = 25
$ws.Range("N33").Value = 23000
$ws.Range("O33").Value = 23000
$ws.Range("P33").Value = 23000
$ws.Range("S33").Value = 1150

# Row 34
$ws.Range("D34").Value = 44435
$ws.Range("M34").Value = 100
$ws.Range("N34").Value = 24000
$ws.Range("O34").Value = 24000
$ws.Range("P34").Value = 24000
$ws.Range("S34").Value = 1200

# Row 35
$ws.Range("D35").Value = 44175
$ws.Range("M35").Value = 25
$ws.Range("N35").Value = 23000
$ws.Range("O35").Value = 23000
$ws.Range("P35").Value = 23000
$ws.Range("S35").Value = 1150

# Row 36
$ws.Range("D36").Value = 44412
$ws.Range("N36").Value = 25000
$ws.Range("O36").Value = 25000
$ws.Range("P36").Value = 25000
$ws.Range("S36").Value = 1250

# Row 37
$ws.Range("D37").Value = 44238
$ws.Range("N37").Value = 25000
$ws.Range("O37").Value = 25000
$ws.Range("P37").Value = 25000
$ws.Range("S37").Value = 1250

# New row 38 (appended observation)
$ws.Range("A38").Value = 10
$ws.Range("B38").Value = "Vega Modelo de Temuco"
$ws.Range("C38").Value = "La Araucanía"
$ws.Range("D38").Value = 44400
$ws.Range("E38").Value = 9
$ws.Range("F38").Value = "Fruta"
$ws.Range("G38").Value = 100108
$ws.Range("H38").Value = "Tropicales y subtropicales"
$ws.Range("I38").Value = 100108007
$ws.Range("J38").Value = "Coco"
$ws.Range("K38").Value = "Sin especificar"
$ws.Range("L38").Value = "Primera"
$ws.Range("M38").Value = 5
$ws.Range("N38").Value = 24000
$ws.Range("O38").Value = 24000
$ws.Range("P38").Value = 24000
$ws.Range("Q38").Value = "$/malla 20 unidades"
$ws.Range("R38").Value = "Perú"
$ws.Range("S38").Value = 1200
$ws.Range("T38").Value = 20

$ws.Range("D38").NumberFormat = $ws.Range("D2").NumberFormat
